$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving it as text, even if it
# looks like a number (e.g. "317.12"), and without leaving a lasting
# text-format override on the cell (format is restored afterwards).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '43.142.12'
$ws.Range("E2").Value = '  +1.94%  '

Set-TextValue $ws.Range("D3") '2.553.70'
$ws.Range("E3").Value = '  +1.25%  '

Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  +0.09%  '

Set-TextValue $ws.Range("D5") '317.12'
$ws.Range("E5").Value = '  +0.25%  '

Set-TextValue $ws.Range("D6") '97.23'
$ws.Range("E6").Value = '  +3.91%  '

$ws.Range("E7").Value = '  +0.64%  '

$ws.Range("E8").Value = '  -0.07%  '

Set-TextValue $ws.Range("D9") '0.542'
$ws.Range("E9").Value = '  +3.15%  '

Set-TextValue $ws.Range("D10") '35.70'
$ws.Range("E10").Value = '  +0.96%  '

Set-TextValue $ws.Range("D11") '0.0811'
$ws.Range("E11").Value = '  +1.02%  '

Set-TextValue $ws.Range("D12") '7.50'
$ws.Range("E12").Value = '  -0.79%  '

$ws.Range("E13").Value = '  -4.60%  '

Set-TextValue $ws.Range("D14") '2.947.33'
$ws.Range("E14").Value = '  +1.27%  '

Set-TextValue $ws.Range("D15") '2.574.41'
$ws.Range("E15").Value = '  +1.89%  '

Set-TextValue $ws.Range("D16") '15.03'
$ws.Range("E16").Value = '  -1.66%  '

Set-TextValue $ws.Range("D17") '0.849'
$ws.Range("E17").Value = '  +1.23%  '

Set-TextValue $ws.Range("D18") '43.185.37'
$ws.Range("E18").Value = '  +1.86%  '

Set-TextValue $ws.Range("D19") '6.83'
$ws.Range("E19").Value = '  +4.52%  '

Set-TextValue $ws.Range("D20") '12.60'
$ws.Range("E20").Value = '  -1.77%  '

Set-TextValue $ws.Range("D21") '0.0₃0964'
$ws.Range("E21").Value = '  +1.15%  '

Set-TextValue $ws.Range("D22") '70.14'
$ws.Range("E22").Value = '  -0.46%  '

Set-TextValue $ws.Range("D23") '253.64'
$ws.Range("E23").Value = '  +1.68%  '

$ws.Range("E24").Value = '  +0.64%  '

Set-TextValue $ws.Range("D25") '2.06'
$ws.Range("E25").Value = '  +2.95%  '

Set-TextValue $ws.Range("D26") '26.76'
$ws.Range("E26").Value = '  +2.34%  '

$ws.Range("E27").Value = '  +0.23%  '

Set-TextValue $ws.Range("D28") '2.43'

Set-TextValue $ws.Range("D29") '40.89'
$ws.Range("E29").Value = '  +5.65%  '

Set-TextValue $ws.Range("D30") '10.26'
$ws.Range("E30").Value = '  +1.64%  '

Set-TextValue $ws.Range("D31") '5.84'
$ws.Range("E31").Value = '  -0.73%  '

Set-TextValue $ws.Range("D32") '156.00'
$ws.Range("E32").Value = '  -0.15%  '

Set-TextValue $ws.Range("D33") '19.34'
$ws.Range("E33").Value = '  -0.07%  '

Set-TextValue $ws.Range("D34") '2.70'
$ws.Range("E34").Value = '  +3.30%  '

$ws.Range("E35").Value = '  +1.35%  '

Set-TextValue $ws.Range("D36") '0.0801'
$ws.Range("E36").Value = '  +2.95%  '

$ws.Range("E37").Value = '  +1.55%  '

$ws.Range("E38").Value = '  +1.60%  '

$ws.Range("E39").Value = '  +5.75%  '

$ws.Range("E40").Value = '  +0.26%  '

Set-TextValue $ws.Range("D41") '22.03'
$ws.Range("E41").Value = '  -6.32%  '

Set-TextValue $ws.Range("D42") '3.88'
$ws.Range("E42").Value = '  +3.08%  '

$ws.Range("E43").Value = '  +2.15%  '

$ws.Range("E44").Value = '  +0.16%  '

$ws.Range("E45").Value = '  -0.86%  '

Set-TextValue $ws.Range("D46") '1.988.21'
$ws.Range("E46").Value = '  -0.95%  '

$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue $ws.Range("D47") '84.63'
$ws.Range("E47").Value = '  +0.64%  '

$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D48") '9.04'
$ws.Range("E48").Value = '  +2.97%  '

Set-TextValue $ws.Range("D49") '2.804.51'
$ws.Range("E49").Value = '  +1.54%  '

Set-TextValue $ws.Range("D50") '104.70'
$ws.Range("E50").Value = '  +3.12%  '

Set-TextValue $ws.Range("D51") '74.43'
$ws.Range("E51").Value = '  +2.55%  '
